$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pitch")

# --- Add the three e-mail addresses + hyperlinks in column A (rows 2-4) ---
$emails = @(
    @{ Row = 2; Mail = "bbb.aaa@ensae.fr" },
    @{ Row = 3; Mail = "ccc.aba@ensae.fr" },
    @{ Row = 4; Mail = "uuu.vvv@ensae.fr" }
)

foreach ($item in $emails) {
    $cell = $ws.Cells.Item($item.Row, 1)
    $cell.Value = $item.Mail
    $ws.Hyperlinks.Add($cell, "mailto:" + $item.Mail)
    # Hyperlinks.Add mints a brand-new cell style (bold hyperlink look); put the
    # alignment back the way it was so the engine re-uses the existing
    # "hyperlink-like" style (s="4") that was already sitting on A2 instead of
    # minting a duplicate one.
    $cell.HorizontalAlignment = -4131
    $cell.VerticalAlignment = -4160
    $cell.WrapText = $true
}

# --- Update the active selection on the "pitch" sheet ---
$ws.Activate()
$ws.Range("A5").Select()
